# Update notes for costs to reflect new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sheet view (scroll position + active selection) ---
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D26").Select()

# --- Rows that previously held the short "benefit is 85% of fee" note
#     now read "Benefit = 85% of fee" and lose their explicit black-font
#     style (s="4" -> default). ---
$plainBenefitRows = @(6, 7, 9, 24, 25, 26, 27, 28, 29, 30, 32, 33, 34, 35, 37, 38, 40)
foreach ($r in $plainBenefitRows) {
    $cell = $ws.Range("E$r")
    $cell.Value = "Benefit = 85% of fee"
    $cell.ClearFormats()
}

# --- Rows 4 and 5 already had no explicit style; just update the text. ---
$ws.Range("E4").Value = "Benefit = 85% of fee"
$ws.Range("E5").Value = "Benefit = 85% of fee"

# --- New note added to the previously-empty E8 cell. ---
$ws.Range("E8").Value = "Emergency department admission - where does this number come from?"

# --- D23 gains a Source value it was previously missing. ---
$ws.Range("D23").Value = "MBS"

# --- Longer, reworded notes (keep their existing styling). ---
$ws.Range("E31").Value = "Benefit = 85% of fee for sum of MBS codes 30473 (`$171.50), 20902 (`$76.70), 72818 (`$91.00)"
$ws.Range("E39").Value = "Benefit = 85% of fee. MBS codes for single lower leg ultrasound (left or right) are 55884 and 55885 (depending on who refers). Values same for knee, ankle etc."

# --- Row 36's note rewritten (no style change, stayed unstyled). ---
$ws.Range("E36").Value = "Benefit = 85% of fee. Note: some oddities about how this is used. I have used this to lump together renal function test costs with ANA and EUC"
